$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# "Added notes from mtg 12/13" - bump the offset formula in column C
# (B+10000 -> B+200000) for every data row, and refresh the view state
# (zoom/selection) to match where the author was working.

for ($r = 2; $r -le 39; $r++) {
    $ws.Cells.Item($r, 3).Formula = "=B$r+200000"
}

# Update sheet view: zoom to 94% and select C14 (author's last-saved cursor)
$excel.ActiveWindow.Zoom = 94
$ws.Range("C14").Select()
